$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Enter the task list text in the same order the shared-string table was
# originally built (B3, B6, B4, B7, B8, B9, B5, B10, B11, B12) so the
# resulting xl/sharedStrings.xml ordering matches the source workbook.
$ws.Range("B3").Value = "Tarefas da semana:"
$ws.Range("B6").Value = "Fazer a hierarquia de ficheiros e adiconá-los"
$ws.Range("B4").Value = "Jogar o jogo para conhecer melhor o projeto"
$ws.Range("B7").Value = "Fazer fork ao projeto e adicionar os membros ao mesmo"
$ws.Range("B8").Value = "Alterar o readme file no repositório git"
$ws.Range("B9").Value = "Meeting semanal"
$ws.Range("B5").Value = "Dar ideias no servidor de discord e discuti-las"
$ws.Range("B10").Value = "Fazer servidor de discord para o trabalho e organiza-lo"
$ws.Range("B11").Value = "Analisar o código dado"
$ws.Range("B12").Value = "Começar a fazer o use case diagram"

# Header row is bold.
$ws.Range("B3").Font.Bold = $true

# Trailing empty, underlined cell left behind near the bottom of the sheet.
$ws.Range("B15").Font.Underline = $true

# F9 was touched (formatting toggled back off) while the sheet was being
# built, which is why the saved used-range/dimension stretches out to
# column F even though the cell itself is empty.
$ws.Range("F9").Font.Bold = $true
$ws.Range("F9").Font.Bold = $false

# Widen the task-description column.
$ws.Columns.Item(2).ColumnWidth = 69.8

# Portrait / A4-ish (PaperSize 9 = A4) page setup for printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the active selection sitting on B8, matching the saved view state.
$ws.Range("B8").Select() | Out-Null
